$wb = $excel.ActiveWorkbook

# --- LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 19:56:05"

$ws.Range("A3").Value = "Total filas: 452"

$ws.Range("A173").Value = "10:52:37"
$ws.Range("C173").Value = "26_HERNANDEZ"
$ws.Range("D173").Value = 89

$ws.Range("A175").Value = "12:01:11"
$ws.Range("C175").Value = "215A_EL PATO"
$ws.Range("D175").Value = 20

$ws.Range("A202").Value = "12:50:41"
$ws.Range("C202").Value = "11_ETCHEVERRY"
$ws.Range("D202").Value = 24

$ws.Range("A203").Value = "12:01:11"
$ws.Range("C203").Value = "215D_EL PATO"
$ws.Range("D203").Value = 73

$ws.Range("A216").Value = "13:18:40"
$ws.Range("C216").Value = "16_SANTA ANA"
$ws.Range("D216").Value = 28

$ws.Range("A217").Value = "12:01:11"
$ws.Range("C217").Value = "17_ROMERO"
$ws.Range("D217").Value = 105

$ws.Range("A244").Value = "14:59:23"
$ws.Range("C244").Value = "11_ETCHEVERRY"
$ws.Range("D244").Value = 1

$ws.Range("A245").Value = "13:18:40"
$ws.Range("C245").Value = "81_EL PELIGRO"
$ws.Range("D245").Value = 102

$ws.Range("A256").Value = "14:59:23"
$ws.Range("C256").Value = "16_SANTA ANA"
$ws.Range("D256").Value = 27

$ws.Range("A257").Value = "13:51:32"
$ws.Range("C257").Value = "26_HERNANDEZ"
$ws.Range("D257").Value = 95

$ws.Range("A314").Value = "15:36:13"
$ws.Range("C314").Value = "215A_EL PATO"
$ws.Range("D314").Value = 89

$ws.Range("A316").Value = "16:53:02"
$ws.Range("C316").Value = "11_ETCHEVERRY"
$ws.Range("D316").Value = 12

$ws.Range("A373").Value = "17:50:30"
$ws.Range("C373").Value = "15_ABASTO"
$ws.Range("D373").Value = 50

$ws.Range("A374").Value = "17:36:40"
$ws.Range("C374").Value = "14_ABASTO"
$ws.Range("D374").Value = 64

$ws.Range("A387").Value = "18:13:12"
$ws.Range("C387").Value = "215A_EL PATO"
$ws.Range("D387").Value = 46

$ws.Range("A389").Value = "18:56:36"
$ws.Range("C389").Value = "11_ETCHEVERRY"
$ws.Range("D389").Value = 3

$ws.Range("A400").Value = "19:16:50"
$ws.Range("C400").Value = "10_OLMOS"
$ws.Range("D400").Value = 4

$ws.Range("A401").Value = "18:35:28"
$ws.Range("C401").Value = "16_SANTA ANA"
$ws.Range("D401").Value = 45

$ws.Range("A402").Value = "17:50:30"
$ws.Range("C402").Value = "14_ABASTO"
$ws.Range("D402").Value = 90

$ws.Range("A403").Value = "18:48:53"
$ws.Range("C403").Value = "26_HERNANDEZ"
$ws.Range("D403").Value = 32

$ws.Range("A416").Value = "18:35:28"
$ws.Range("C416").Value = "16_P MOR-SANTA ANA"
$ws.Range("D416").Value = 75

$ws.Range("A417").Value = "18:13:12"
$ws.Range("C417").Value = "11X44_ETCHEVERRY"
$ws.Range("D417").Value = 97

$ws.Range("A421").Value = "19:56:05"
$ws.Range("B421").Value = "19:56"
$ws.Range("C421").Value = "16_SANTA ANA"
$ws.Range("D421").Value = 0

$ws.Range("A422").Value = "18:13:12"
$ws.Range("B422").Value = "19:59"
$ws.Range("C422").Value = "17_ROMERO"
$ws.Range("D422").Value = 106

$ws.Range("A423").Value = "19:42:01"
$ws.Range("B423").Value = "20:01"
$ws.Range("C423").Value = "16_SANTA ANA"
$ws.Range("D423").Value = 19

$ws.Range("A424").Value = "18:56:36"
$ws.Range("B424").Value = "20:09"
$ws.Range("C424").Value = "15_ABASTO"
$ws.Range("D424").Value = 73

$ws.Range("A425").Value = "18:35:28"
$ws.Range("B425").Value = "20:10"
$ws.Range("D425").Value = 95

$ws.Range("A426").Value = "18:13:12"
$ws.Range("C426").Value = "16_P MOR-167 Y 521"
$ws.Range("D426").Value = 118

$ws.Range("A427").Value = "19:42:01"
$ws.Range("B427").Value = "20:11"
$ws.Range("C427").Value = "10_OLMOS"
$ws.Range("D427").Value = 29

$ws.Range("B428").Value = "20:12"
$ws.Range("C428").Value = "23_HERNANDEZ"
$ws.Range("D428").Value = 56

$ws.Range("A429").Value = "19:16:50"
$ws.Range("C429").Value = "14_ABASTO"
$ws.Range("D429").Value = 57

$ws.Range("A430").Value = "19:42:01"
$ws.Range("B430").Value = "20:13"
$ws.Range("C430").Value = "23_HERNANDEZ"
$ws.Range("D430").Value = 31

$ws.Range("B431").Value = "20:20"
$ws.Range("C431").Value = "26_HERNANDEZ"
$ws.Range("D431").Value = 92

$ws.Range("A433").Value = "18:48:53"
$ws.Range("B433").Value = "20:21"
$ws.Range("D433").Value = 93

$ws.Range("B434").Value = "20:22"
$ws.Range("C434").Value = "11_ETCHEVERRY"
$ws.Range("D434").Value = 107

$ws.Range("A435").Value = "18:35:28"
$ws.Range("B435").Value = "20:23"
$ws.Range("D435").Value = 108

$ws.Range("A436").Value = "19:56:05"
$ws.Range("B436").Value = "20:23"
$ws.Range("C436").Value = "16_SANTA ANA"
$ws.Range("D436").Value = 27

$ws.Range("A437").Value = "19:42:01"
$ws.Range("B437").Value = "20:24"
$ws.Range("C437").Value = "215A_EL PATO"
$ws.Range("D437").Value = 42

$ws.Range("A438").Value = "18:48:53"
$ws.Range("B438").Value = "20:30"
$ws.Range("C438").Value = "225_GOMEZ"
$ws.Range("D438").Value = 102

$ws.Range("A439").Value = "18:35:28"
$ws.Range("B439").Value = "20:31"
$ws.Range("C439").Value = "225_GOMEZ"
$ws.Range("D439").Value = 116

$ws.Range("A440").Value = "19:42:01"
$ws.Range("B440").Value = "20:32"
$ws.Range("C440").Value = "14_ABASTO"
$ws.Range("D440").Value = 50

$ws.Range("A441").Value = "19:56:05"
$ws.Range("B441").Value = "20:34"
$ws.Range("D441").Value = 38

$ws.Range("B442").Value = "20:44"
$ws.Range("C442").Value = "11_ETCHEVERRY"
$ws.Range("D442").Value = 62

$ws.Range("A443").Value = "19:16:50"
$ws.Range("B443").Value = "20:49"
$ws.Range("C443").Value = "11_ETCHEVERRY"
$ws.Range("D443").Value = 93

$ws.Range("B444").Value = "20:50"
$ws.Range("C444").Value = "14_ABASTO"
$ws.Range("D444").Value = 94

$ws.Range("A445").Value = "19:42:01"
$ws.Range("B445").Value = "20:52"
$ws.Range("C445").Value = "15_ABASTO"
$ws.Range("D445").Value = 70

$ws.Range("B446").Value = "20:53"
$ws.Range("C446").Value = "23_HERNANDEZ"
$ws.Range("D446").Value = 71

$ws.Range("B447").Value = "20:55"
$ws.Range("C447").Value = "10_OLMOS"
$ws.Range("D447").Value = 99

$ws.Range("B448").Value = "20:56"
$ws.Range("C448").Value = "27_EL RETIRO"
$ws.Range("D448").Value = 100

$ws.Range("B449").Value = "20:57"
$ws.Range("C449").Value = "27_EL RETIRO"
$ws.Range("D449").Value = 75

$ws.Range("A450").Value = "19:16:50"
$ws.Range("B450").Value = "21:04"
$ws.Range("C450").Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Range("D450").Value = 108

$ws.Range("A451").Value = "19:16:50"
$ws.Range("B451").Value = "21:07"
$ws.Range("C451").Value = "215B_EL PATO"
$ws.Range("D451").Value = 111

$ws.Range("B452").Value = "21:08"
$ws.Range("C452").Value = "215B_EL PATO"
$ws.Range("D452").Value = 86

$ws.Range("B453").Value = "21:21"
$ws.Range("C453").Value = "26_HERNANDEZ"
$ws.Range("D453").Value = 99

$ws.Range("A454").Value = "19:42:01"
$ws.Range("B454").Value = "21:23"
$ws.Range("C454").Value = "10_OLMOS"
$ws.Range("D454").Value = 101
$ws.Range("E454").Value = "LP1912"

$ws.Range("A455").Value = "19:42:01"
$ws.Range("B455").Value = "21:38"
$ws.Range("C455").Value = "14_ABASTO"
$ws.Range("D455").Value = 116
$ws.Range("E455").Value = "LP1912"

$ws.Range("A456").Value = "19:42:01"
$ws.Range("B456").Value = "21:38"
$ws.Range("C456").Value = "17_ROMERO"
$ws.Range("D456").Value = 116
$ws.Range("E456").Value = "LP1912"

$ws.Range("A457").Value = "19:56:05"
$ws.Range("B457").Value = "21:47"
$ws.Range("C457").Value = "215A_EL PATO"
$ws.Range("D457").Value = 111
$ws.Range("E457").Value = "LP1912"

# --- LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 19:56:05"

$ws.Range("A3").Value = "Total filas: 48"

$ws.Range("A53").Value = "19:56:05"
$ws.Range("B53").Value = "21:47"
$ws.Range("C53").Value = "215A_EL PATO"
$ws.Range("D53").Value = 111
$ws.Range("E53").Value = "LP1912"

# --- 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 19:56:05"
